$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.611.97'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = '1.587.23'
$ws.Range('E3').Value = '  -2.27%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'211.02"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('D6').Value = "'0.509"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.62%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').Value = "'19.51"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.00%  '
$ws.Range('D11').Value = "'0.0834"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').Value = '1.809.26'
$ws.Range('E12').Value = '  -2.27%  '
$ws.Range('D13').Value = '1.590.66'
$ws.Range('E13').Value = '  -2.05%  '
$ws.Range('E14').Value = '  -2.91%  '
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = "'64.77"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '26.592.10'
$ws.Range('E17').Value = '  -1.97%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -2.51%  '
$ws.Range('D19').Value = "'209.06"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.13%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  -3.13%  '
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('D25').Value = "'145.44"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('E28').Value = '  -2.72%  '
$ws.Range('D29').Value = "'15.26"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('E33').Value = '  +22.63%  '
$ws.Range('D34').Value = "'2.90"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('D35').Value = '1.309.57'
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').Value = "'1.49"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.12%  '
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('E41').Value = '  +2.88%  '
$ws.Range('D42').Value = "'0.790"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('D43').Value = "'2.16"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('D44').Value = "'62.57"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.62%  '
$ws.Range('D45').Value = '1.722.70'
$ws.Range('E45').Value = '  -2.07%  '
$ws.Range('D46').Value = "'89.34"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('D47').Value = "'1.61"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.10%  '
$ws.Range('D48').Value = "'0.840"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.70%  '
$ws.Range('E49').Value = '  -1.70%  '
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('D51').Value = "'7.52"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.15%  '
